$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 34.349525
$ws.Range("H2").Value = 103.048575
$ws.Range("I2").Value = 0.9193849879779717
$ws.Range("J2").Value = 0.9193849879779719
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 127.984071
$ws.Range("N2").Value = 383.952213
$ws.Range("O2").Value = 0.4594739502473105
$ws.Range("P2").Value = 0.4594739502473105
$ws.Range("Q2").Value = 4396.192046416275
$ws.Range("R2").Value = 39565.72841774648
$ws.Range("S2").Value = 0.4224334522243148
$ws.Range("T2").Value = 0.4224334522243148

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 34.349525
$ws.Range("H3").Value = 103.048575
$ws.Range("I3").Value = 0.9193849879779717
$ws.Range("J3").Value = 0.9193849879779719
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 60.45343933333334
$ws.Range("N3").Value = 181.360318
$ws.Range("O3").Value = 0.2170331070069088
$ws.Range("P3").Value = 0.2170331070069088
$ws.Range("Q3").Value = 2076.546925716317
$ws.Range("R3").Value = 18688.92233144685
$ws.Range("S3").Value = 0.1995369804763687
$ws.Range("T3").Value = 0.1995369804763687

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 34.349525
$ws.Range("H4").Value = 103.048575
$ws.Range("I4").Value = 0.9193849879779717
$ws.Range("J4").Value = 0.9193849879779719
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 59.37981533333334
$ws.Range("N4").Value = 178.139446
$ws.Range("O4").Value = 0.2131787034353868
$ws.Range("P4").Value = 0.2131787034353868
$ws.Range("Q4").Value = 2039.668451287717
$ws.Range("R4").Value = 18357.01606158945
$ws.Range("S4").Value = 0.1959932996951027
$ws.Range("T4").Value = 0.1959932996951027

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 34.349525
$ws.Range("H5").Value = 103.048575
$ws.Range("I5").Value = 0.9193849879779717
$ws.Range("J5").Value = 0.9193849879779719
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 28.25780433333334
$ws.Range("N5").Value = 84.77341300000001
$ws.Range("O5").Value = 0.1014479761497213
$ws.Range("P5").Value = 0.1014479761497213
$ws.Range("Q5").Value = 970.6421563929417
$ws.Range("R5").Value = 8735.779407536476
$ws.Range("S5").Value = 0.09326974633280108
$ws.Range("T5").Value = 0.09326974633280109

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 34.349525
$ws.Range("H6").Value = 103.048575
$ws.Range("I6").Value = 0.9193849879779717
$ws.Range("J6").Value = 0.9193849879779719
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 2.469651333333333
$ws.Range("N6").Value = 7.408954
$ws.Range("O6").Value = 0.008866263160672582
$ws.Range("P6").Value = 0.008866263160672582
$ws.Range("Q6").Value = 84.83135021561667
$ws.Range("R6").Value = 763.48215194055
$ws.Range("S6").Value = 0.008151509249384495
$ws.Range("T6").Value = 0.008151509249384497

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 1.096366
$ws.Range("H7").Value = 3.289098
$ws.Range("I7").Value = 0.02934487279604178
$ws.Range("J7").Value = 0.02934487279604178
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 127.984071
$ws.Range("N7").Value = 383.952213
$ws.Range("O7").Value = 0.4594739502473105
$ws.Range("P7").Value = 0.4594739502473105
$ws.Range("Q7").Value = 140.317383985986
$ws.Range("R7").Value = 1262.856455873874
$ws.Range("S7").Value = 0.01348320462310216
$ws.Range("T7").Value = 0.01348320462310216

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 1.096366
$ws.Range("H8").Value = 3.289098
$ws.Range("I8").Value = 0.02934487279604178
$ws.Range("J8").Value = 0.02934487279604178
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 60.45343933333334
$ws.Range("N8").Value = 181.360318
$ws.Range("O8").Value = 0.2170331070069088
$ws.Range("P8").Value = 0.2170331070069088
$ws.Range("Q8").Value = 66.27909546812933
$ws.Range("R8").Value = 596.511859213164
$ws.Range("S8").Value = 0.006368808917647461
$ws.Range("T8").Value = 0.006368808917647462

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 1.096366
$ws.Range("H9").Value = 3.289098
$ws.Range("I9").Value = 0.02934487279604178
$ws.Range("J9").Value = 0.02934487279604178
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 59.37981533333334
$ws.Range("N9").Value = 178.139446
$ws.Range("O9").Value = 0.2131787034353868
$ws.Range("P9").Value = 0.2131787034353868
$ws.Range("Q9").Value = 65.10201061774534
$ws.Range("R9").Value = 585.918095559708
$ws.Range("S9").Value = 0.006255701935136541
$ws.Range("T9").Value = 0.006255701935136541

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 1.096366
$ws.Range("H10").Value = 3.289098
$ws.Range("I10").Value = 0.02934487279604178
$ws.Range("J10").Value = 0.02934487279604178
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 28.25780433333334
$ws.Range("N10").Value = 84.77341300000001
$ws.Range("O10").Value = 0.1014479761497213
$ws.Range("P10").Value = 0.1014479761497213
$ws.Range("Q10").Value = 30.98089590571934
$ws.Range("R10").Value = 278.828063151474
$ws.Range("S10").Value = 0.002976977955529452
$ws.Range("T10").Value = 0.002976977955529452

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 1.096366
$ws.Range("H11").Value = 3.289098
$ws.Range("I11").Value = 0.02934487279604178
$ws.Range("J11").Value = 0.02934487279604178
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 2.469651333333333
$ws.Range("N11").Value = 7.408954
$ws.Range("O11").Value = 0.008866263160672582
$ws.Range("P11").Value = 0.008866263160672582
$ws.Range("Q11").Value = 2.707641753721333
$ws.Range("R11").Value = 24.368775783492
$ws.Range("S11").Value = 0.0002601793646261683
$ws.Range("T11").Value = 0.0002601793646261683

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 0.9739683333333332
$ws.Range("H12").Value = 2.921905
$ws.Range("I12").Value = 0.02606882815504994
$ws.Range("J12").Value = 0.02606882815504994
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 127.984071
$ws.Range("N12").Value = 383.952213
$ws.Range("O12").Value = 0.4594739502473105
$ws.Range("P12").Value = 0.4594739502473105
$ws.Range("Q12").Value = 124.652432325085
$ws.Range("R12").Value = 1121.871890925765
$ws.Range("S12").Value = 0.0119779474507191
$ws.Range("T12").Value = 0.0119779474507191

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 0.9739683333333332
$ws.Range("H13").Value = 2.921905
$ws.Range("I13").Value = 0.02606882815504994
$ws.Range("J13").Value = 0.02606882815504994
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 60.45343933333334
$ws.Range("N13").Value = 181.360318
$ws.Range("O13").Value = 0.2170331070069088
$ws.Range("P13").Value = 0.2170331070069088
$ws.Range("Q13").Value = 58.87973555175444
$ws.Range("R13").Value = 529.91761996579
$ws.Range("S13").Value = 0.005657798770519668
$ws.Range("T13").Value = 0.005657798770519669

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 0.9739683333333332
$ws.Range("H14").Value = 2.921905
$ws.Range("I14").Value = 0.02606882815504994
$ws.Range("J14").Value = 0.02606882815504994
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 59.37981533333334
$ws.Range("N14").Value = 178.139446
$ws.Range("O14").Value = 0.2131787034353868
$ws.Range("P14").Value = 0.2131787034353868
$ws.Range("Q14").Value = 57.83405977384778
$ws.Range("R14").Value = 520.5065379646301
$ws.Range("S14").Value = 0.005557318986173453
$ws.Range("T14").Value = 0.005557318986173453

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 0.9739683333333332
$ws.Range("H15").Value = 2.921905
$ws.Range("I15").Value = 0.02606882815504994
$ws.Range("J15").Value = 0.02606882815504994
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 28.25780433333334
$ws.Range("N15").Value = 84.77341300000001
$ws.Range("O15").Value = 0.1014479761497213
$ws.Range("P15").Value = 0.1014479761497213
$ws.Range("Q15").Value = 27.52220659019611
$ws.Range("R15").Value = 247.699859311765
$ws.Range("S15").Value = 0.002644629856924689
$ws.Range("T15").Value = 0.00264462985692469

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 0.9739683333333332
$ws.Range("H16").Value = 2.921905
$ws.Range("I16").Value = 0.02606882815504994
$ws.Range("J16").Value = 0.02606882815504994
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 2.469651333333333
$ws.Range("N16").Value = 7.408954
$ws.Range("O16").Value = 0.008866263160672582
$ws.Range("P16").Value = 0.008866263160672582
$ws.Range("Q16").Value = 2.405362193041111
$ws.Range("R16").Value = 21.64825973737
$ws.Range("S16").Value = 0.0002311330907130234
$ws.Range("T16").Value = 0.0002311330907130235

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 0.9415566666666667
$ws.Range("H17").Value = 2.82467
$ws.Range("I17").Value = 0.02520131107093657
$ws.Range("J17").Value = 0.02520131107093657
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 127.984071
$ws.Range("N17").Value = 383.952213
$ws.Range("O17").Value = 0.4594739502473105
$ws.Range("P17").Value = 0.4594739502473105
$ws.Range("Q17").Value = 120.50425527719
$ws.Range("R17").Value = 1084.53829749471
$ws.Range("S17").Value = 0.01157934594917451
$ws.Range("T17").Value = 0.01157934594917451

$ws.Range("E18").Value = 3
$ws.Range("G18").Value = 0.9415566666666667
$ws.Range("H18").Value = 2.82467
$ws.Range("I18").Value = 0.02520131107093657
$ws.Range("J18").Value = 0.02520131107093657
$ws.Range("K18").Value = 3
$ws.Range("M18").Value = 60.45343933333334
$ws.Range("N18").Value = 181.360318
$ws.Range("O18").Value = 0.2170331070069088
$ws.Range("P18").Value = 0.2170331070069088
$ws.Range("Q18").Value = 56.9203388272289
$ws.Range("R18").Value = 512.2830494450601
$ws.Range("S18").Value = 0.00546951884237297
$ws.Range("T18").Value = 0.005469518842372971

$ws.Range("E19").Value = 3
$ws.Range("G19").Value = 0.9415566666666667
$ws.Range("H19").Value = 2.82467
$ws.Range("I19").Value = 0.02520131107093657
$ws.Range("J19").Value = 0.02520131107093657
$ws.Range("K19").Value = 3
$ws.Range("M19").Value = 59.37981533333334
$ws.Range("N19").Value = 178.139446
$ws.Range("O19").Value = 0.2131787034353868
$ws.Range("P19").Value = 0.2131787034353868
$ws.Range("Q19").Value = 55.90946099253556
$ws.Range("R19").Value = 503.1851489328201
$ws.Range("S19").Value = 0.005372382818974118
$ws.Range("T19").Value = 0.005372382818974118

$ws.Range("E20").Value = 3
$ws.Range("G20").Value = 0.9415566666666667
$ws.Range("H20").Value = 2.82467
$ws.Range("I20").Value = 0.02520131107093657
$ws.Range("J20").Value = 0.02520131107093657
$ws.Range("K20").Value = 3
$ws.Range("M20").Value = 28.25780433333334
$ws.Range("N20").Value = 84.77341300000001
$ws.Range("O20").Value = 0.1014479761497213
$ws.Range("P20").Value = 0.1014479761497213
$ws.Range("Q20").Value = 26.60632405541223
$ws.Range("R20").Value = 239.45691649871
$ws.Range("S20").Value = 0.00255662200446608
$ws.Range("T20").Value = 0.002556622004466081

$ws.Range("E21").Value = 3
$ws.Range("G21").Value = 0.9415566666666667
$ws.Range("H21").Value = 2.82467
$ws.Range("I21").Value = 0.02520131107093657
$ws.Range("J21").Value = 0.02520131107093657
$ws.Range("K21").Value = 3
$ws.Range("M21").Value = 2.469651333333333
$ws.Range("N21").Value = 7.408954
$ws.Range("O21").Value = 0.008866263160672582
$ws.Range("P21").Value = 0.008866263160672582
$ws.Range("Q21").Value = 2.325316677242222
$ws.Range("R21").Value = 20.92785009518
$ws.Range("S21").Value = 0.000223441455948895
$ws.Range("T21").Value = 0.000223441455948895
